# Fix Training Data Issue (#48)
# Data in column BF ("Date") was off by one day due to the way NBA stats
# were shown; update values from "5-5-2012-13" to "2013-05-05".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 58).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 31 }

$dateCol = $ws.Range($ws.Cells.Item(2, 58), $ws.Cells.Item($lastRow, 58))

# Force text interpretation so the corrected value isn't silently
# re-parsed into a date serial by the COM value setter.
$dateCol.NumberFormat = "@"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # Column BF = 58
    if ($cell.Value2 -eq "5-5-2012-13") {
        $cell.Value = "2013-05-05"
    }
}

# Restore the cells' original (default) formatting now that the text
# value is safely stored.
$dateCol.ClearFormats()
